# fix typo array notation
#
# Slide 4, shape "TextBox 5" contains a Java code sample whose 4th
# line reads:
#     int a1[] = { 5, 1, 7, 10 };
# The array brackets belong on the type, not the variable name, so the
# line should read:
#     int[] a1 = { 5, 1, 7, 10 };
#
# The line is split across two runs: "int" (bold) and
# " a1[] = { 5, 1, 7, 10 };" (regular). We edit each run's text via
# Characters(start, length) sub-ranges of the shape's TextRange so the
# existing run-level formatting (bold/Courier New/etc.) is preserved.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

# Paragraph 4 is "4\t\t\tint a1[] = { 5, 1, 7, 10 };"
$para = $tr.Paragraphs(4, 1)

# "int" starts right after the leading "4\t\t\t" (4 characters).
$intStart = $para.Start + 4
$intRange = $tr.Characters($intStart, 3)
$intRange.Text = "int[]"

# Re-locate the paragraph (its length just grew) and grab everything
# after the new "int[]" run, i.e. " a1[] = { 5, 1, 7, 10 };".
$restStart = $intStart + 5
$para2 = $tr.Paragraphs(4, 1)
$restLen = ($para2.Start + $para2.Length) - $restStart
$restRange = $tr.Characters($restStart, $restLen)
$restRange.Text = " a1 = { 5, 1, 7, 10 };"
